$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 4 of the table corresponds to "Test Case #3" (currently empty requirement/description/output cells)
$t.Cell(4, 2).Range.Text = "Order list cannot be empty when clicking " + [char]8220 + "remove" + [char]8221 + " button in current order window."
$t.Cell(4, 3).Range.Text = "Clicking on remove button with empty order list in current order window."
$t.Cell(4, 4).Range.Text = "Error message."

# Row 5 of the table corresponds to "Test Case #4" (currently empty requirement/description/output cells)
$t.Cell(5, 2).Range.Text = "Order list cannot be empty when clicking " + [char]8220 + "complete order" + [char]8221 + " button in current order window."
$t.Cell(5, 3).Range.Text = "Clicking on add to order button with empty order list in current order window."
$t.Cell(5, 4).Range.Text = "Error message."
